# Update workbook/sheet for new "through" date (2022-08-04 -> 2022-08-05)
# and add the new day's data for August (column I = "Total"/current year column).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new date
$ws.Name = "Through 2022-08-05"

# Update header label in I1 ("2022 (through 08-04)" -> "2022 (through 08-05)")
$ws.Range("I1").Value = "2022 (through 08-05)"

# Update the August row (row 9) current-year total
$ws.Range("I9").Value = 25

# Update the Total row (row 14) current-year total
$ws.Range("I14").Value = 995
